# Update "想去人数" (want-to-go count) / "最低票价" (min ticket price) figures
# across the three sheets that carry this data: 展览 (Exhibition), 演出
# (Performance) and 全部类型 (All types, a combined view of 演出+展览+本地生活).
# 本地生活 is untouched by this commit.

$wb = $excel.ActiveWorkbook

# ---- 展览 ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 269
$ws.Range("F4").Value  = 281
$ws.Range("F5").Value  = 2887
$ws.Range("F8").Value  = 2246
$ws.Range("F9").Value  = 1443
$ws.Range("F10").Value = 1443
$ws.Range("F12").Value = 451
$ws.Range("G12").Value = 80
$ws.Range("F14").Value = 2583
$ws.Range("F16").Value = 1405
$ws.Range("F17").Value = 4971
$ws.Range("F19").Value = 5367
$ws.Range("F20").Value = 5367
$ws.Range("F21").Value = 1927
$ws.Range("F23").Value = 3356
$ws.Range("F25").Value = 1627
$ws.Range("F26").Value = 272
$ws.Range("F27").Value = 847
$ws.Range("F28").Value = 138
$ws.Range("F29").Value = 5
$ws.Range("F30").Value = 328
$ws.Range("F31").Value = 1046
$ws.Range("F32").Value = 2149
$ws.Range("F34").Value = 132
$ws.Range("F35").Value = 306
$ws.Range("F36").Value = 792
$ws.Range("F37").Value = 168
$ws.Range("F38").Value = 375
$ws.Range("F39").Value = 446

# ---- 演出 ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 60

# ---- 全部类型 (演出 rows 2-6, 展览 rows 7-44, 本地生活 rows 45-51) -------
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value  = 60
$ws.Range("F7").Value  = 269
$ws.Range("F9").Value  = 281
$ws.Range("F10").Value = 2887
$ws.Range("F12").Value = 2246
$ws.Range("F13").Value = 1443
$ws.Range("F14").Value = 1443
$ws.Range("F17").Value = 451
$ws.Range("G17").Value = 80
$ws.Range("F20").Value = 2583
$ws.Range("F21").Value = 1405
$ws.Range("F26").Value = 4971
$ws.Range("F28").Value = 5367
$ws.Range("F29").Value = 5367
$ws.Range("F30").Value = 1927
$ws.Range("F32").Value = 3356
$ws.Range("F36").Value = 1627
$ws.Range("F38").Value = 272
$ws.Range("F39").Value = 847
$ws.Range("F40").Value = 138
$ws.Range("F41").Value = 5
$ws.Range("F42").Value = 328
$ws.Range("F44").Value = 2149
$ws.Range("F46").Value = 132
$ws.Range("F47").Value = 306
$ws.Range("F48").Value = 792
$ws.Range("F49").Value = 168
$ws.Range("F50").Value = 375
$ws.Range("F51").Value = 446
